# "treat wire transfers only optionally as exchanges, address some further feedback"
#
# On the "currency_movements" sheet, split the old "currency" column into a
# separate "fees" column (new) and a "currency" column (shifted right), and
# reduce the two withdrawal rows' "amount" by the newly-broken-out fee.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_movements")

# Insert a new, blank column in front of the old "currency" column (D) so the
# old D ("currency") and E ("comment") columns shift right to E and F.
$ws.Columns("D:D").Insert()

# The inserted column inherits formatting from the column to its left; strip
# that back off so the new "fees" cells start out unstyled like the rest of
# the plain numeric cells on this sheet.
$ws.Range("D2:D5").ClearFormats()

# New header + fee values for the (previously non-existent) "fees" column.
$ws.Range("D1").Value = "fees"
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0

# The wire-transfer "amount" on the two outgoing rows now excludes the fee
# that used to be folded into it (-100 -> -98, i.e. a $2 fee broken out).
$ws.Range("C4").Value = -98
$ws.Range("C5").Value = -98
